{"js": "// Fix PICE-related issues in the forecast: corrected dates (June 8 -> 9,\n// June 10 -> 11, June 16 -> 17) and corrected sea-ice figures (concentration\n// -> coverage, date range, and percentages).\n\nconst replacements = [\n  // Narrative paragraph: first significant pulse (15%) date.\n  [\"is expected on June 8, which\", \"is expected on June 9, which\"],\n  // Narrative paragraph: half-way point (50%) date.\n  [\"be reached on June 16, which\", \"be reached on June 17, which\"],\n  // Narrative paragraph: sea-ice concentration -> coverage figures.\n  [\n    \"the Spring (April 20 \\u2013 May 31) sea ice concentration was 0 %, significantly lower than the long-term average of 0.55 %, as measured\",\n    \"the Spring (March 19 \\u2013 May 31) sea ice coverage was 12 %, significantly lower than the long-term average of 55 %, as measured\",\n  ],\n  // \"Run Timing Forecast by Date\" paragraph: \"concentrations\" -> \"coverage\".\n  [\"the sea ice concentrations and the air\", \"the sea ice coverage and the air\"],\n  // \"Run Timing Forecast by Date\" paragraph: the three predicted dates.\n  [\n    \"are June 8 (15%), June 10 (25%), and June 16 (50%) in District\",\n    \"are June 9 (15%), June 11 (25%), and June 17 (50%) in District\",\n  ],\n];\n\nfor (const [searchText, replaceText] of replacements) {\n  const results = context.document.body.search(searchText, { matchCase: true });\n  results.load(\"items\");\n  await context.sync();\n  for (const item of results.items) {\n    item.insertText(replaceText, \"Replace\");\n  }\n  await context.sync();\n}\n\n// Table cells on the \"Run Timing Forecast by Date\" table give the same\n// three dates as standalone cell text (15% / 25% / 50% rows).\nconst tableDateReplacements = [\n  [\"June 8\", \"June 9\"],\n  [\"June 10\", \"June 11\"],\n  [\"June 16\", \"June 17\"],\n];\n\nfor (const [searchText, replaceText] of tableDateReplacements) {\n  const results = context.document.body.search(searchText, {\n    matchCase: true,\n    matchWholeWord: true,\n  });\n  results.load(\"items\");\n  await context.sync();\n  for (const item of results.items) {\n    item.insertText(replaceText, \"Replace\");\n  }\n  await context.sync();\n}\n", "ps1": "# Fix PICE-related issues in the forecast: corrected dates (June 8 -> 9,\n# June 10 -> 11, June 16 -> 17) and corrected sea-ice figures (concentration\n# -> coverage, date range, and percentages).\n\n$d = $word.ActiveDocument\n\n$replacements = @(\n    # Narrative paragraph: first significant pulse (15%) date.\n    @(\"is expected on June 8, which\", \"is expected on June 9, which\"),\n    # Narrative paragraph: half-way point (50%) date.\n    @(\"be reached on June 16, which\", \"be reached on June 17, which\"),\n    # Narrative paragraph: sea-ice concentration -> coverage figures.\n    @(\"the Spring (April 20 \" + [char]0x2013 + \" May 31) sea ice concentration was 0 %, significantly lower than the long-term average of 0.55 %, as measured\",\n      \"the Spring (March 19 \" + [char]0x2013 + \" May 31) sea ice coverage was 12 %, significantly lower than the long-term average of 55 %, as measured\"),\n    # \"Run Timing Forecast by Date\" paragraph: \"concentrations\" -> \"coverage\".\n    @(\"the sea ice concentrations and the air\", \"the sea ice coverage and the air\"),\n    # \"Run Timing Forecast by Date\" paragraph: the three predicted dates.\n    @(\"are June 8 (15%), June 10 (25%), and June 16 (50%) in District\",\n      \"are June 9 (15%), June 11 (25%), and June 17 (50%) in District\")\n)\n\nforeach ($pair in $replacements) {\n    $findText = $pair[0]\n    $replaceText = $pair[1]\n    $r = $d.Content\n    $r.Find.Execute($findText, $false, $false, $false, $false, $false, $true, 1, $false, $replaceText, 2)\n}\n\n# Table cells on the \"Run Timing Forecast by Date\" table give the same\n# three dates as standalone cell text (15% / 25% / 50% rows).\n$tableDateReplacements = @(\n    @(\"June 8\", \"June 9\"),\n    @(\"June 10\", \"June 11\"),\n    @(\"June 16\", \"June 17\")\n)\n\nforeach ($pair in $tableDateReplacements) {\n    $findText = $pair[0]\n    $replaceText = $pair[1]\n    $r = $d.Content\n    $r.Find.Execute($findText, $false, $true, $false, $false, $false, $true, 1, $false, $replaceText, 2)\n}\n"}
